$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.770179333333333
$ws.Range("H2").Value = 5.310538
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8369776666666665
$ws.Range("N2").Value = 2.510933
$ws.Range("O2").Value = 0.0694586718035551
$ws.Range("P2").Value = 0.06945867180355511
$ws.Range("Q2").Value = 1.481600567994889
$ws.Range("R2").Value = 13.334405111954
$ws.Range("S2").Value = 0.0694586718035551
$ws.Range("T2").Value = 0.06945867180355511

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.770179333333333
$ws.Range("H3").Value = 5.310538
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.427350333333333
$ws.Range("N3").Value = 7.282051
$ws.Range("O3").Value = 0.2014397000898671
$ws.Range("P3").Value = 0.2014397000898671
$ws.Range("Q3").Value = 4.296845394826445
$ws.Range("R3").Value = 38.671608553438
$ws.Range("S3").Value = 0.2014397000898671
$ws.Range("T3").Value = 0.2014397000898671

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.770179333333333
$ws.Range("H4").Value = 5.310538
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.785681666666667
$ws.Range("N4").Value = 26.357045
$ws.Range("O4").Value = 0.7291016281065776
$ws.Range("P4").Value = 0.7291016281065776
$ws.Range("Q4").Value = 15.55223211557889
$ws.Range("R4").Value = 139.97008904021
$ws.Range("S4").Value = 0.7291016281065776
$ws.Range("T4").Value = 0.7291016281065776
